$d = $word.ActiveDocument

$replacements = @(
    @("2026-01-17 Saturday", "2026-01-18 Sunday"),
    @("135÷3=45, 0", "379÷6=63, 1"),
    @("284÷7=40, 4", "845÷6=140, 5"),
    @("515÷8=64, 3", "165÷6=27, 3"),
    @("434÷5=86, 4", "640÷8=80, 0"),
    @("307÷4=76, 3", "991÷7=141, 4"),
    @("869÷9=96, 5", "143÷7=20, 3"),
    @("930÷6=155, 0", "751÷3=250, 1"),
    @("950÷9=105, 5", "474÷3=158, 0"),
    @("362÷5=72, 2", "284÷4=71, 0"),
    @("779÷9=86, 5", "809÷9=89, 8"),
    @("158÷6=26, 2", "369÷4=92, 1"),
    @("249÷8=31, 1", "823÷7=117, 4"),
    @("908÷5=181, 3", "797÷2=398, 1"),
    @("440÷3=146, 2", "646÷5=129, 1"),
    @("410÷8=51, 2", "118÷8=14, 6"),
    @("436÷4=109, 0", "710÷8=88, 6"),
    @("167÷4=41, 3", "223÷2=111, 1"),
    @("400÷9=44, 4", "225÷2=112, 1"),
    @("805÷5=161, 0", "572÷2=286, 0"),
    @("672÷2=336, 0", "103÷3=34, 1"),
    @("190÷8=23, 6", "870÷5=174, 0"),
    @("985÷9=109, 4", "214÷8=26, 6"),
    @("198÷9=22, 0", "463÷6=77, 1"),
    @("315÷2=157, 1", "130÷4=32, 2"),
    @("279÷8=34, 7", "260÷8=32, 4")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Write-Output "Applied $($replacements.Count) replacements"
